# Update the "Notes:" narrative on the About sheet to reflect the revised
# rationale for the 3% discount rate (now referencing the Social Cost of
# Carbon central estimate, rather than only the health-damages discount
# rate), and drop the trailing note about the EU's 1%-based SCoC estimate.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("A16").Value = "We choose to use a 3% discount rate here, for consistency with the 3% rate used for the central estimate"
$ws.Range("A17").Value = "of Social Cost of Carbon (in the SCoC variable), as well as the discount rate built into the health"
$ws.Range("A18").Value = "damages values in the SCoHIbP Social Cost of Health Impacts by Pollutant variable."
$ws.Range("A19").ClearContents()
